# Append a new row (17) to the sheet, recording a submission at
# 2025-05-01T11:48:55.200Z (mirrors the other "raw ISO timestamp" rows,
# e.g. row 16 / row 6-8, that this form-export sheet already contains).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 17

$ws.Cells.Item($row, 1).Value = "2025-05-01T11:48:55.200Z"
$ws.Cells.Item($row, 2).Value = "UNDP"
$ws.Cells.Item($row, 3).Value = "C4"
$ws.Cells.Item($row, 4).Value = "الرحلة 3"
$ws.Cells.Item($row, 5).Value = "الصمود"
$ws.Cells.Item($row, 6).Value = "يامن "
# "43434" reads as a plain number unless entered as text (leading
# apostrophe) - the column holds free-text IDs elsewhere in the sheet
# (e.g. G16 = "421123"), so keep it text here too.
$ws.Cells.Item($row, 7).Value = "'43434"
# Column H ("الوقت") is blank for this submission, same as most other rows.
$ws.Cells.Item($row, 8).Value = ""
